# Select the "Misc" worksheet (sheet4.xml / rId4)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")

# Insert a new row before the current row 324, shifting existing rows down.
$ws.Rows.Item(324).Insert()

# Populate the newly inserted row with the new card entry.
$ws.Range("A324").Value = "Mickey Moniak 2017 Topps Pro Debut"
$ws.Range("B324").Value = "https://blowoutbuzz.files.wordpress.com/2016/12/2017-topps-pro-debut-moniak-auto.jpg"

# Update the active selection to match the author's final cursor position.
$ws.Activate()
$ws.Range("B324").Select()
